# Re-run with updated costs: refresh computed funding/per-pupil/change metrics
# for rows 4-12 (data rows) on the active sheet, per the latest model output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("I4").Value = 3411944290.748064
$ws.Range("J4").Value = 3345363634.868221
$ws.Range("K4").Value = 3318503399.207038
$ws.Range("L4").Value = 3336501063.008941
$ws.Range("M4").Value = 3315039799.618937
$ws.Range("N4").Value = 3512740323.367759
$ws.Range("O4").Value = 0.727326968973747
$ws.Range("P4").Value = 0.7285202863961814
$ws.Range("Q4").Value = 0.704653937947494
$ws.Range("R4").Value = 0.7673031026252983
$ws.Range("S4").Value = 0.5363961813842482
$ws.Range("T4").Value = 6236.42370491668
$ws.Range("U4").Value = 6014.807286148247
$ws.Range("V4").Value = 5900.379784587087
$ws.Range("W4").Value = 5875.818128841674
$ws.Range("X4").Value = 5787.182286405236
$ws.Range("Y4").Value = 6110.165616980592
$ws.Range("Z4").Value = -221.6164187684335
$ws.Range("AA4").Value = -336.0439203295937
$ws.Range("AB4").Value = -360.6055760750069
$ws.Range("AC4").Value = -449.241418511444
$ws.Range("AD4").Value = -126.2580879360885
$ws.Range("AE4").Value = -0.03553581816349571
$ws.Range("AF4").Value = -0.05388407462832634
$ws.Range("AG4").Value = -0.05782249461188982
$ws.Range("AH4").Value = -0.07203510213029152
$ws.Range("AI4").Value = -0.02024527099346196
# Row 5
$ws.Range("I5").Value = 4492589279.743372
$ws.Range("J5").Value = 4439928328.463231
$ws.Range("K5").Value = 4399244281.204008
$ws.Range("L5").Value = 4421954997.706984
$ws.Range("M5").Value = 4331657750.37428
$ws.Range("N5").Value = 4606016896.876914
$ws.Range("O5").Value = 0.7111785533636823
$ws.Range("P5").Value = 0.7516439049064239
$ws.Range("Q5").Value = 0.7546788062721295
$ws.Range("R5").Value = 0.813859382903389
$ws.Range("S5").Value = 0.6469398077895802
$ws.Range("T5").Value = 6264.299056983244
$ws.Range("U5").Value = 6085.624096257529
$ws.Range("V5").Value = 5938.528748366973
$ws.Range("W5").Value = 5883.384853664924
$ws.Range("X5").Value = 5696.102030842227
$ws.Range("Y5").Value = 6020.110830378724
$ws.Range("Z5").Value = -178.6749607257143
$ws.Range("AA5").Value = -325.770308616271
$ws.Range("AB5").Value = -380.9142033183198
$ws.Range("AC5").Value = -568.197026141017
$ws.Range("AD5").Value = -244.1882266045195
$ws.Range("AE5").Value = -0.02852273799516847
$ws.Range("AF5").Value = -0.05200427145206499
$ws.Range("AG5").Value = -0.06080715493518607
$ws.Range("AH5").Value = -0.09070400709998172
$ws.Range("AI5").Value = -0.03898093376182388
# Row 6
$ws.Range("I6").Value = 8276601906.599
$ws.Range("J6").Value = 8109547338.001504
$ws.Range("K6").Value = 8045840242.740464
$ws.Range("L6").Value = 8088271392.502736
$ws.Range("M6").Value = 7854145149.792853
$ws.Range("N6").Value = 8114336299.011919
$ws.Range("O6").Value = 0.7645403377110694
$ws.Range("P6").Value = 0.7692307692307693
$ws.Range("Q6").Value = 0.74812382739212
$ws.Range("R6").Value = 0.8208255159474672
$ws.Range("S6").Value = 0.6880863039399625
$ws.Range("T6").Value = 7730.73472749087
$ws.Range("U6").Value = 7463.234925104884
$ws.Range("V6").Value = 7329.201877009309
$ws.Range("W6").Value = 7311.583782083367
$ws.Range("X6").Value = 7079.106835657934
$ws.Range("Y6").Value = 7348.521186357456
$ws.Range("Z6").Value = -267.4998023859862
$ws.Range("AA6").Value = -401.5328504815607
$ws.Range("AB6").Value = -419.1509454075031
$ws.Range("AC6").Value = -651.6278918329363
$ws.Range("AD6").Value = -382.2135411334139
$ws.Range("AE6").Value = -0.03460211892082443
$ws.Range("AF6").Value = -0.05193980451220115
$ws.Range("AG6").Value = -0.05421877223609572
$ws.Range("AH6").Value = -0.08429055126102258
$ws.Range("AI6").Value = -0.04944077821920911
# Row 7
$ws.Range("I7").Value = 1829119704.047315
$ws.Range("J7").Value = 1767572366.328891
$ws.Range("K7").Value = 1773451732.905504
$ws.Range("L7").Value = 1783856268.46569
$ws.Range("M7").Value = 1763982872.839947
$ws.Range("N7").Value = 1846979964.771495
$ws.Range("P7").Value = 0.7456037514654161
$ws.Range("Q7").Value = 0.7268464243845252
$ws.Range("R7").Value = 0.7690504103165299
$ws.Range("S7").Value = 0.6400937866354045
$ws.Range("T7").Value = 6715.557316402897
$ws.Range("U7").Value = 6425.888513413438
$ws.Range("V7").Value = 6367.755454679464
$ws.Range("W7").Value = 6331.703906754919
$ws.Range("X7").Value = 6200.574623233917
$ws.Range("Y7").Value = 6476.881681733364
$ws.Range("Z7").Value = -289.6688029894594
$ws.Range("AA7").Value = -347.8018617234329
$ws.Range("AB7").Value = -383.853409647978
$ws.Range("AC7").Value = -514.9826931689804
$ws.Range("AD7").Value = -238.6756346695338
$ws.Range("AE7").Value = -0.04313399310611754
$ws.Range("AF7").Value = -0.05179046880798999
$ws.Range("AG7").Value = -0.0571588315850432
$ws.Range("AH7").Value = -0.07668502685713419
$ws.Range("AI7").Value = -0.03554070398395126
# Row 8
$ws.Range("C8").Value = 923251
$ws.Range("D8").Value = 936427
$ws.Range("E8").Value = 950363
$ws.Range("F8").Value = 961402.5
$ws.Range("G8").Value = 970645
$ws.Range("H8").Value = 974400
$ws.Range("I8").Value = 5959152943.871764
$ws.Range("J8").Value = 5810259650.273874
$ws.Range("K8").Value = 5820587323.240251
$ws.Range("L8").Value = 5850284223.438849
$ws.Range("M8").Value = 5801456586.886761
$ws.Range("N8").Value = 6043277356.649564
$ws.Range("O8").Value = 0.770893371757925
$ws.Range("P8").Value = 0.7921469740634006
$ws.Range("Q8").Value = 0.7881844380403458
$ws.Range("R8").Value = 0.8332132564841499
$ws.Range("S8").Value = 0.6884005763688761
$ws.Range("T8").Value = 6454.531805404776
$ws.Range("U8").Value = 6204.711793096391
$ws.Range("V8").Value = 6124.59378494349
$ws.Range("W8").Value = 6085.156033439532
$ws.Range("X8").Value = 5976.908743038661
$ws.Range("Y8").Value = 6202.049832357927
$ws.Range("Z8").Value = -249.8200123083843
$ws.Range("AA8").Value = -329.9380204612853
$ws.Range("AB8").Value = -369.375771965244
$ws.Range("AC8").Value = -477.6230623661149
$ws.Range("AD8").Value = -252.4819730468489
$ws.Range("AE8").Value = -0.03870459079606592
$ws.Range("AF8").Value = -0.05111726619504886
$ws.Range("AG8").Value = -0.05722735329244844
$ws.Range("AH8").Value = -0.07399809572030802
$ws.Range("AI8").Value = -0.03911700812062469
# Row 9
$ws.Range("B9").Value = 2794
$ws.Range("C9").Value = 1051493.5
$ws.Range("D9").Value = 1070056
$ws.Range("E9").Value = 1085270.5
$ws.Range("F9").Value = 1099211.5
$ws.Range("G9").Value = 1112866
$ws.Range("H9").Value = 1121019
$ws.Range("I9").Value = 6433931311.054861
$ws.Range("J9").Value = 6348008285.547463
$ws.Range("K9").Value = 6344550385.639488
$ws.Range("L9").Value = 6430165491.548478
$ws.Range("M9").Value = 6363876968.513569
$ws.Range("N9").Value = 6730323471.86068
$ws.Range("O9").Value = 0.7072297780959198
$ws.Range("P9").Value = 0.7559055118110236
$ws.Range("Q9").Value = 0.7168933428775949
$ws.Range("R9").Value = 0.7759484609878311
$ws.Range("S9").Value = 0.5801717967072297
$ws.Range("T9").Value = 6118.850293468158
$ws.Range("U9").Value = 5932.40754273371
$ws.Range("V9").Value = 5846.054403615954
$ws.Range("W9").Value = 5849.798234050934
$ws.Range("X9").Value = 5718.457539823814
$ws.Range("Y9").Value = 6003.755040602059
$ws.Range("Z9").Value = -186.4427507344481
$ws.Range("AA9").Value = -272.7958898522047
$ws.Range("AB9").Value = -269.052059417224
$ws.Range("AC9").Value = -400.3927536443443
$ws.Range("AD9").Value = -115.0952528660991
$ws.Range("AE9").Value = -0.03047022590722226
$ws.Range("AF9").Value = -0.04458286716761362
$ws.Range("AG9").Value = -0.04397101522559488
$ws.Range("AH9").Value = -0.06543594538859066
$ws.Range("AI9").Value = -0.01880994751399012
# Row 10
$ws.Range("B10").Value = 1706
$ws.Range("C10").Value = 570762
$ws.Range("D10").Value = 578735
$ws.Range("E10").Value = 584523
$ws.Range("F10").Value = 590081
$ws.Range("G10").Value = 594464.5
$ws.Range("H10").Value = 596000.5
$ws.Range("I10").Value = 3506623273.280505
$ws.Range("J10").Value = 3455335077.422535
$ws.Range("K10").Value = 3420868913.955131
$ws.Range("L10").Value = 3448441759.298356
$ws.Range("M10").Value = 3392540506.063832
$ws.Range("N10").Value = 3610553223.380827
$ws.Range("O10").Value = 0.6729191090269636
$ws.Range("P10").Value = 0.7004689331770223
$ws.Range("Q10").Value = 0.6547479484173505
$ws.Range("R10").Value = 0.7391559202813599
$ws.Range("S10").Value = 0.5228604923798359
$ws.Range("T10").Value = 6143.757421272799
$ws.Range("U10").Value = 5970.496129355466
$ws.Range("V10").Value = 5852.411135156583
$ws.Range("W10").Value = 5844.014227365999
$ws.Range("X10").Value = 5706.884946138638
$ws.Range("Y10").Value = 6057.970124825109
$ws.Range("Z10").Value = -173.2612919173325
$ws.Range("AA10").Value = -291.3462861162161
$ws.Range("AB10").Value = -299.7431939068001
$ws.Range("AC10").Value = -436.8724751341606
$ws.Range("AD10").Value = -85.78729644768919
$ws.Range("AE10").Value = -0.02820119351024741
$ws.Range("AF10").Value = -0.0474215152290075
$ws.Range("AG10").Value = -0.04878825340156456
$ws.Range("AH10").Value = -0.07110835359831869
$ws.Range("AI10").Value = -0.01396332741762396
# Row 11
$ws.Range("I11").Value = 4693698414.146843
$ws.Range("J11").Value = 4598180429.977851
$ws.Range("K11").Value = 4603332041.951568
$ws.Range("L11").Value = 4635506918.00151
$ws.Range("M11").Value = 4582083727.028511
$ws.Range("N11").Value = 4779957404.706257
$ws.Range("O11").Value = 0.7306101344364012
$ws.Range("P11").Value = 0.7285418821096173
$ws.Range("Q11").Value = 0.7119958634953464
$ws.Range("R11").Value = 0.7699069286452948
$ws.Range("S11").Value = 0.6240951396070321
$ws.Range("T11").Value = 6521.194778189193
$ws.Range("U11").Value = 6280.907882008597
$ws.Range("V11").Value = 6210.065356731347
$ws.Range("W11").Value = 6172.569249720713
$ws.Range("X11").Value = 6047.103365079621
$ws.Range("Y11").Value = 6288.606738481429
$ws.Range("Z11").Value = -240.2868961805962
$ws.Range("AA11").Value = -311.1294214578456
$ws.Range("AB11").Value = -348.6255284684803
$ws.Range("AC11").Value = -474.0914131095724
$ws.Range("AD11").Value = -232.5880397077644
$ws.Range("AE11").Value = -0.03684706627445944
$ws.Range("AF11").Value = -0.04771049355839674
$ws.Range("AG11").Value = -0.05346037656082503
$ws.Range("AH11").Value = -0.07270008475980816
$ws.Range("AI11").Value = -0.03566647640792442
# Row 12
$ws.Range("B12").Value = 1727
$ws.Range("C12").Value = 626429.5
$ws.Range("D12").Value = 635047
$ws.Range("E12").Value = 641895.5
$ws.Range("F12").Value = 648990
$ws.Range("G12").Value = 652427
$ws.Range("H12").Value = 652097
$ws.Range("I12").Value = 4047638095.340007
$ws.Range("J12").Value = 3902603647.954151
$ws.Range("K12").Value = 3891547801.611442
$ws.Range("L12").Value = 3918131835.529537
$ws.Range("M12").Value = 3891307093.331212
$ws.Range("N12").Value = 4070846790.876235
$ws.Range("O12").Value = 0.7672264041690793
$ws.Range("P12").Value = 0.7990735379270412
$ws.Range("Q12").Value = 0.7730167921250723
$ws.Range("R12").Value = 0.8048639258830341
$ws.Range("S12").Value = 0.6282570932252461
$ws.Range("T12").Value = 6461.442341620258
$ws.Range("U12").Value = 6145.37766173866
$ws.Range("V12").Value = 6062.587760175046
$ws.Range("W12").Value = 6037.276129878021
$ws.Range("X12").Value = 5964.356308569712
$ws.Range("Y12").Value = 6242.7013019171
$ws.Range("Z12").Value = -316.0646798815978
$ws.Range("AA12").Value = -398.8545814452118
$ws.Range("AB12").Value = -424.1662117422366
$ws.Range("AC12").Value = -497.0860330505457
$ws.Range("AD12").Value = -218.7410397031581
$ws.Range("AE12").Value = -0.0489154995388138
$ws.Range("AF12").Value = -0.06172841300092691
$ws.Range("AG12").Value = -0.06564574739142115
$ws.Range("AH12").Value = -0.07693112570991345
$ws.Range("AI12").Value = -0.03385328354540529
